$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing cell text (org/TEN TO CHUC + CAP columns for rows 2 & 3)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "BAN PHONG TRÀO TÌNH NGUYỆN"
$ws.Range("D2").Value = "Ban Chuyên môn"

$ws.Range("C3").Value = "KTX 135"
$ws.Range("D3").Value = "Khoa/Viện/KTX"

# ---------------------------------------------------------------------------
# 2. Add the new 3rd data row (row 4) - copy formats from row 3 so the new
#    row visually matches the rest of the table, then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "org.003@gmail.com"
$ws.Range("C4").Value = "CLB BẠN TRAI"
$ws.Range("D4").Value = "CLB/Đội/Nhóm"
$ws.Range("E4").Value = "DontKnowWTODO"

$h = $ws.Hyperlinks.Add($ws.Range("B4"), "mailto:org.003@gmail.com")

$ws.Rows.Item(4).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3. Column C is now a lot wider to fit the new text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 39.140625

# ---------------------------------------------------------------------------
# 4. Misc view / print bits.
# ---------------------------------------------------------------------------
$ws.Range("E10").Select()
$ws.PageSetup.Orientation = 1
